$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- LTSD Parameters table (K2:N5) ---
# K2 currently holds "Test"; repurpose it as the table title.
$ws.Range("K2").Value2 = "LTSD Parameters"

# Row 3: side headers
$ws.Range("K3").Value2 = "Right"
$ws.Range("M3").Value2 = "Left"

# Row 4: column headers
$ws.Range("K4").Value2 = "Threshols"
$ws.Range("L4").Value2 = "Win"
$ws.Range("M4").Value2 = "Threshold"
$ws.Range("N4").Value2 = "Win"

# Row 5: values. These must stay TEXT (not get auto-converted to numbers),
# so format the cells as Text first, write the values, then drop the
# number-format again so the cells end up back on the default style.
$ws.Range("K5:N5").NumberFormat = "@"
$ws.Range("K5").Value2 = "5.0"
$ws.Range("L5").Value2 = "200.0"
$ws.Range("M5").Value2 = "6.0"
$ws.Range("N5").Value2 = "300.0"
$ws.Range("K5:N5").Style = "Normal"

# Put the selection where the author's cursor ended up.
[void]$ws.Range("N5").Select()
